$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CISF" (sheet1.xml): updated objective-function data + new constraint
# column ("prebuild") and a new "Superstorage" site row.
# ---------------------------------------------------------------------------
$cisf = $wb.Worksheets.Item("CISF")

# New header for column D
$cisf.Range("D1").Value = "prebuild"

# capacity (col B) changed from 5000 to 420 for the three existing sites
$cisf.Range("B2").Value = 420
$cisf.Range("B3").Value = 420
$cisf.Range("B4").Value = 420

# new "prebuild" flag column, set to 1 for the existing sites
$cisf.Range("D2").Value = 1
$cisf.Range("D3").Value = 1
$cisf.Range("D4").Value = 1

# new row: Superstorage site
$cisf.Range("A5").Value = "Superstorage"
$cisf.Range("B5").Value = 500
$cisf.Range("C5").Value = 999
$cisf.Range("D5").Value = 0

[void]$cisf.Activate()
$cisf.Range("C6").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "General" (sheet2.xml): no data change, only the stored selection.
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("General")
[void]$general.Activate()
$general.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Reactors" (sheet4.xml): objective function simplified to only
# name/snf columns, and a new "Gorleben" row is now part of the reactor list.
# ---------------------------------------------------------------------------
$reactors = $wb.Worksheets.Item("Reactors")

# drop the now-unused capacity/costs columns
$reactors.Range("C1:D3").ClearContents()

# new row for Gorleben
$reactors.Range("A4").Value = "Gorleben"
$reactors.Range("B4").Value = 2000

[void]$reactors.Activate()
$reactors.Range("A5").Select() | Out-Null
